$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9714348912239075
$ws.Range("B1").Value = 0.8539865612983704
$ws.Range("C1").Value = 0.659881055355072
$ws.Range("D1").Value = 0.6383389830589294
$ws.Range("E1").Value = 0.6881824135780334
